# Update latest output (run 122)
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Schedule": collapse the 3 schedule rows into a single row, and
# update its values. Row 2 is rewritten in place; rows 3 and 4 are removed.
# ---------------------------------------------------------------------------
$wsSchedule = $wb.Worksheets.Item("Schedule")

$wsSchedule.Range("A2").Value = 46042
$wsSchedule.Range("B2").Value = 46042.66666666666
$wsSchedule.Range("C2").Value = 16
$wsSchedule.Range("D2").Value = 60.48
$wsSchedule.Range("E2").Value = 807.3803789999996
$wsSchedule.Range("F2").Value = 13.34954330357142

$wsSchedule.Rows("3:4").Delete()

# ---------------------------------------------------------------------------
# Sheet "Detailed": refreshed historical/forecast price series + pump status.
# ---------------------------------------------------------------------------
$wsDetailed = $wb.Worksheets.Item("Detailed")

$wsDetailed.Range("E2").Value = "ON"

$wsDetailed.Range("E11").Value = "ON"
$wsDetailed.Range("E12").Value = "ON"
$wsDetailed.Range("E13").Value = "ON"
$wsDetailed.Range("E14").Value = "ON"
$wsDetailed.Range("E15").Value = "ON"

$wsDetailed.Range("B16").Value = 56.98
$wsDetailed.Range("E16").Value = "ON"

$wsDetailed.Range("B17").Value = 8.61645
$wsDetailed.Range("E17").Value = "ON"

$wsDetailed.Range("B18").Value = 0.01928

$wsDetailed.Range("B19").Value = 8.682
$wsDetailed.Range("C19").Value = "historical"

$wsDetailed.Range("B20").Value = -5.51
$wsDetailed.Range("C20").Value = "historical"

$wsDetailed.Range("B21").Value = -6.19659
$wsDetailed.Range("B22").Value = -6.32
$wsDetailed.Range("B23").Value = -6.36387
$wsDetailed.Range("B24").Value = -7.22044
$wsDetailed.Range("B25").Value = -7.25296
$wsDetailed.Range("B26").Value = -6.09234
$wsDetailed.Range("B27").Value = -7.2065
$wsDetailed.Range("B28").Value = -8.84238
$wsDetailed.Range("B29").Value = -7.79956
$wsDetailed.Range("B30").Value = -7.37432
$wsDetailed.Range("B31").Value = -7.37879
$wsDetailed.Range("B32").Value = -6.36144
$wsDetailed.Range("B33").Value = -5.51
$wsDetailed.Range("B34").Value = 0.7
$wsDetailed.Range("B35").Value = -1.11816
$wsDetailed.Range("B36").Value = 0
$wsDetailed.Range("B37").Value = 9.916370000000001
$wsDetailed.Range("B38").Value = 9.979480000000001
$wsDetailed.Range("B39").Value = 33.15828

$wsDetailed.Range("B41").Value = 57.3

$wsDetailed.Range("B42").Value = 59.77564
$wsDetailed.Range("E42").Value = "OFF"

$wsDetailed.Range("B43").Value = 59.66166
$wsDetailed.Range("E43").Value = "OFF"

$wsDetailed.Range("B44").Value = 57.91165
$wsDetailed.Range("E44").Value = "OFF"

$wsDetailed.Range("B45").Value = 56.98
$wsDetailed.Range("E45").Value = "OFF"

$wsDetailed.Range("E46").Value = "OFF"
$wsDetailed.Range("E47").Value = "OFF"

$wsDetailed.Range("B48").Value = 57.3
$wsDetailed.Range("E48").Value = "OFF"

$wsDetailed.Range("B49").Value = 57.06003
$wsDetailed.Range("E49").Value = "OFF"
